$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated latitude (V) / longitude (W) coordinates for rows 60-131
# as part of the automated "Actualizar" data refresh.
$coords = @(
    @(60, 15.5341, -88.019300000000001),
    @(61, 15.537100000000001, -88.014899999999997),
    @(62, 15.542899999999999, -88.023099999999999),
    @(63, 15.549099999999999, -88.035899999999998),
    @(64, 15.5505, -88.004800000000003),
    @(65, 15.6119, -87.956299999999999),
    @(66, 15.316800000000001, -87.990300000000005),
    @(67, 15.4361, -87.921099999999996),
    @(68, 15.436500000000001, -87.924300000000002),
    @(69, 15.438499999999999, -87.926900000000003),
    @(70, 15.4405, -87.929100000000005),
    @(71, 13.309900000000001, -87.179100000000005),
    @(72, 13.3125, -87.175299999999993),
    @(73, 13.095800000000001, -87.057100000000005),
    @(74, 14.0412, -87.232399999999998),
    @(75, 14.0457, -87.211399999999998),
    @(76, 14.0543, -87.2226),
    @(77, 14.055099999999999, -87.222300000000004),
    @(78, 14.055400000000001, -87.222099999999998),
    @(79, 14.0558, -87.229299999999995),
    @(80, 14.0588, -87.189599999999999),
    @(81, 14.058999999999999, -87.189400000000006),
    @(82, 14.059100000000001, -87.220600000000005),
    @(83, 14.0601, -83.400400000000005),
    @(84, 14.064, -87.209599999999995),
    @(85, 14.0655, -87.1785),
    @(86, 14.0655, -87.179299999999998),
    @(87, 14.0655, -87.179599999999994),
    @(88, 14.0665, -87.210700000000003),
    @(89, 14.0692, -87.185000000000002),
    @(90, 14.081799999999999, -87.209599999999995),
    @(91, 14.083399999999999, -87.174999999999997),
    @(92, 14.0837, -87.210099999999997),
    @(93, 14.0844, -87.182699999999997),
    @(94, 14.085000000000001, -87.186999999999998),
    @(95, 14.085900000000001, -87.174400000000006),
    @(96, 14.0869, -87.184899999999999),
    @(97, 14.088200000000001, -87.183400000000006),
    @(98, 14.089600000000001, -87.188500000000005),
    @(99, 14.09, -87.206500000000005),
    @(100, 14.090199999999999, -87.197199999999995),
    @(101, 14.090400000000001, -87.196799999999996),
    @(102, 14.092599999999999, -87.2393),
    @(103, 14.092700000000001, -87.194100000000006),
    @(104, 14.0928, -87.218500000000006),
    @(105, 14.0947, -87.1952),
    @(106, 14.095700000000001, -87.182699999999997),
    @(107, 14.095700000000001, -87.211200000000005),
    @(108, 14.097, -87.207499999999996),
    @(109, 14.097, -87.222499999999997),
    @(110, 14.099600000000001, -87.194500000000005),
    @(111, 14.0998, -87.189599999999999),
    @(112, 14.0999, -87.187100000000001),
    @(113, 14.1004, -87.208500000000001),
    @(114, 14.1004, -87.183999999999997),
    @(115, 14.1006, -87.182199999999995),
    @(116, 14.1007, -87.180599999999998),
    @(117, 14.1007, -87.183400000000006),
    @(118, 14.101000000000001, -87.182699999999997),
    @(119, 14.1014, -87.207800000000006),
    @(120, 14.102600000000001, -87.184799999999996),
    @(121, 14.1028, -87.178899999999999),
    @(122, 14.1029, -87.196799999999996),
    @(123, 14.1043, -87.198599999999999),
    @(124, 14.1045, -87.199200000000005),
    @(125, 14.1046, -87.236400000000003),
    @(126, 14.1052, -87.205399999999997),
    @(127, 14.106299999999999, -87.204499999999996),
    @(128, 14.1065, -87.205399999999997),
    @(129, 14.1067, -87.206800000000001),
    @(130, 14.1068, -87.205799999999996),
    @(131, 14.107900000000001, -87.205299999999994)
)

foreach ($entry in $coords) {
    $r = $entry[0]
    $ws.Cells.Item($r, 22).Value2 = $entry[1]
    $ws.Cells.Item($r, 23).Value2 = $entry[2]
}

# Restore the frozen-pane scroll position / active selection to the last
# row touched by this refresh (row 132), matching the author's session view.
$win = $excel.ActiveWindow
$win.ScrollRow = 126
$win.ScrollColumn = 19
$ws.Range("A132:XFD132").Select()
